# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) FernandaJourney: two copy text tweaks
# ---------------------------------------------------------------------------
$fj = $wb.Worksheets.Item("FernandaJourney")
$fj.Range("B8").Value = "fuckkkk"
$fj.Range("B11").Value = "I can't resist you anymore"

# ---------------------------------------------------------------------------
# 2) Split "cumcontrol" into "cumcontrol1" (edited copy of the original) and
#    "cumcontrol2" (a second variant, seeded from the "dickpic" sheet's
#    layout/style and populated with new delay/sync/edge copy).
# ---------------------------------------------------------------------------
$cumcontrol1 = $wb.Worksheets.Item("cumcontrol")
$cumcontrol1.Name = "cumcontrol1"

# Rewrite the message/note copy for cumcontrol1 (names in column A stay the same)
$cumcontrol1.Range("B2").Value = "patience... what's coming is worth every second of waiting"

$cumcontrol1.Range("B3").Value = "hold it for me... I have years of experience and this next one is my best work"
$cumcontrol1.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol1.Range("B4").Value = "I want to feel you let go while I do the same... watch this first"
$cumcontrol1.Range("C4").Value = "SYNC variant. Send PPV."

$cumcontrol1.Range("B5").Value = "now we go together amor... I've been holding back too. open this"
$cumcontrol1.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol1.Range("B6").Value = "a man who can wait gets rewarded... trust me on that"

$cumcontrol1.Range("B7").Value = "I can tell you're close... not yet amor, I know what I'm doing"
$cumcontrol1.Range("C7").Value = "CONTROL."

# Seed "cumcontrol2" from the "dickpic" sheet so it inherits the same column
# widths / header & row styling, inserting it directly before "dickpic".
$dickpic = $wb.Worksheets.Item("dickpic")
$dickpic.Copy($dickpic, $null)
$cumcontrol2 = $wb.Worksheets.Item("dickpic (2)")
$cumcontrol2.Name = "cumcontrol2"

$cumcontrol2.Range("A2").Value = "delay2"
$cumcontrol2.Range("B2").Value = "save it for this last one amor, I promise you it's going to be worth it"
$cumcontrol2.Range("C2").Value = "DELAY variant."

$cumcontrol2.Range("A3").Value = "delay1"
$cumcontrol2.Range("B3").Value = "one more for you before we're done... this is the one I'm most proud of"
$cumcontrol2.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol2.Range("A4").Value = "sync2"
$cumcontrol2.Range("B4").Value = "I'm ready when you are... but see this first"
$cumcontrol2.Range("C4").Value = "SYNC variant."

$cumcontrol2.Range("A5").Value = "sync1"
$cumcontrol2.Range("B5").Value = "okay... let's both let go right now. open this"
$cumcontrol2.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol2.Range("A6").Value = "edge2"
$cumcontrol2.Range("B6").Value = "not yet... a little more anticipation makes it so much better, trust me"
$cumcontrol2.Range("C6").Value = "EDGE variant."

$cumcontrol2.Range("A7").Value = "edge1"
$cumcontrol2.Range("B7").Value = "slow down for me amor... I know exactly when to let you go"
$cumcontrol2.Range("C7").Value = "CONTROL."

# The original "dickpic" sheet is left completely untouched, it simply now
# sits after "cumcontrol2" (and before "boosters").
